# ENG SOFWTARE ajustes 12/09/2024 2024.2
#
# Applies the wording/formatting tweaks to the "Diagrama(s) de Classe(s)"
# slides (19-24) plus one textbox resize on slide 22.
#
# Helper: find $old inside the shape's TextRange and overwrite just that
# span via TextRange.Characters(start, length) so the untouched
# surrounding runs (and their bold/color formatting) are left alone.
function Set-SubText {
    param(
        $TextRange,
        [string]$OldText,
        [string]$NewText
    )
    $full = $TextRange.Text
    $idx = $full.IndexOf($OldText)
    $span = $TextRange.Characters($idx + 1, $OldText.Length)
    $span.Text = $NewText
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 19 - "Diagrama de Classe (Instanciado = Objeto)" textbox
# ---------------------------------------------------------------------
$s19 = $p.Slides.Item(19)
$tr19 = $s19.Shapes.Item(4).TextFrame.TextRange
Set-SubText $tr19 "Diagramas de Classes (" "Diagrama de Classe ("
Set-SubText $tr19 "Classes, representação de um item do mundo real." "Classe, representação de um item do mundo real."

# ---------------------------------------------------------------------
# Slide 20 - "Diagrama de Classe (Visibiliadade)" textbox
# ---------------------------------------------------------------------
$s20 = $p.Slides.Item(20)
$tr20 = $s20.Shapes.Item(4).TextFrame.TextRange
Set-SubText $tr20 "Diagramas de Classes (" "Diagrama de Classe ("

# ---------------------------------------------------------------------
# Slide 21 - "Diagrama de Classe (Exemplo)" textbox
# ---------------------------------------------------------------------
$s21 = $p.Slides.Item(21)
$tr21 = $s21.Shapes.Item(4).TextFrame.TextRange
Set-SubText $tr21 "Diagramas de Classes (" "Diagrama de Classe ("
Set-SubText $tr21 "Relacionamentos entre as Classes (conexões)" "Associação entre Classes (conexão) - Relacionamentos"

# ---------------------------------------------------------------------
# Slide 22 - "Diagramas de Classe (Multiplicidade)" textbox
# ---------------------------------------------------------------------
$s22 = $p.Slides.Item(22)
$shape22 = $s22.Shapes.Item(4)

# grow the textbox (only the height changes, position/width stay put)
$shape22.Height = 5209475 / 12700

$tr22 = $shape22.TextFrame.TextRange
Set-SubText $tr22 "Diagramas de Classes (" "Diagramas de Classe ("
Set-SubText $tr22 "Cardinalidade" "Multiplicidade"

# ---------------------------------------------------------------------
# Slide 23 - "Diagrama de Classe (Associação)" textbox
# ---------------------------------------------------------------------
$s23 = $p.Slides.Item(23)
$tr23 = $s23.Shapes.Item(4).TextFrame.TextRange
Set-SubText $tr23 "Diagramas de Classes (" "Diagrama de Classe ("
Set-SubText $tr23 "Relacionamentos" "Associação"

# ---------------------------------------------------------------------
# Slide 24 - "Diagrama de Classe (Associação)" textbox
# here the leading run is split into "Diagrama de Classe " + "("
# ---------------------------------------------------------------------
$s24 = $p.Slides.Item(24)
$tr24 = $s24.Shapes.Item(4).TextFrame.TextRange
Set-SubText $tr24 "Diagramas de Classes " "Diagrama de Classe "
Set-SubText $tr24 "Relacionamentos" "Associação"
